$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7024
$ws1.Range("F4").Value = 462
$ws1.Range("F7").Value = 149
$ws1.Range("F12").Value = 200
$ws1.Range("F17").Value = 3628
$ws1.Range("G17").Value = 65
$ws1.Range("F23").Value = 2255
$ws1.Range("F25").Value = 256
$ws1.Range("F32").Value = 248
$ws1.Range("F33").Value = 94

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7024
$ws4.Range("F4").Value = 462
$ws4.Range("F8").Value = 149
$ws4.Range("F13").Value = 200
$ws4.Range("F18").Value = 3628
$ws4.Range("G18").Value = 65
$ws4.Range("F24").Value = 2255
$ws4.Range("F26").Value = 256
$ws4.Range("F33").Value = 248
$ws4.Range("F34").Value = 94

$wb.Save()
